$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: column D prices are stored as plain text in this sheet (e.g. "1.015"),
# even though many of them look like numbers. Assigning such a string via
# .Value would normally get auto-converted to a number by Excel, so for those
# cells we prefix the value with a literal leading apostrophe (forces text
# entry, same as typing '1.015 into a cell) and then reset the cell's style
# back to "Normal" so no stray formatting/quote-prefix is left behind.

$ws.Range("D2").Value = '27.046.20'
$ws.Range("E2").Value = '  +0.24%  '
$ws.Range("D3").Value = '1.848.37'
$ws.Range("E3").Value = '  +0.38%  '
$ws.Range("E4").Value = '  +0.71%  '
$ws.Range("D5").Value = '''1.015'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.64%  '
$ws.Range("D6").Value = '''309.83'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.12%  '
$ws.Range("D7").Value = '''0.4759'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.96%  '
$ws.Range("D8").Value = '''0.3684'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.79%  '
$ws.Range("D9").Value = '''0.07241'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.61%  '
$ws.Range("D10").Value = '''0.9322'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.37%  '
$ws.Range("D11").Value = '''19.91'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.00%  '
$ws.Range("D12").Value = '''0.07793'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.45%  '
$ws.Range("D13").Value = '1.860.11'
$ws.Range("E13").Value = '  +1.85%  '
$ws.Range("D14").Value = '''5.397'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.44%  '
$ws.Range("D15").Value = '''6.483'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.46%  '
$ws.Range("D16").Value = '''88.99'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.93%  '
$ws.Range("D17").Value = '''1.017'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.68%  '
$ws.Range("D18").Value = '''0.000008661'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.97%  '
$ws.Range("E19").Value = '  +0.74%  '
$ws.Range("D20").Value = '27.090.55'
$ws.Range("E20").Value = '  +0.32%  '
$ws.Range("D21").Value = '''14.56'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.88%  '
$ws.Range("D22").Value = '''5.058'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.93%  '
$ws.Range("D24").Value = '''1.943'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.61%  '
$ws.Range("D25").Value = '''152.92'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.31%  '
$ws.Range("D26").Value = '''18.37'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.89%  '
$ws.Range("E27").Value = '  -1.96%  '
$ws.Range("D28").Value = '''114.66'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.58%  '
$ws.Range("D29").Value = '''4.927'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.91%  '
$ws.Range("D30").Value = '''0.08867'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.13%  '
$ws.Range("D31").Value = '''3.323'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.74%  '
$ws.Range("D32").Value = '''1.180'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.95%  '
$ws.Range("B33").Value = 'ImmutableX'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D33").Value = '''0.7389'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.98%  '
$ws.Range("B34").Value = 'Filecoin'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D34").Value = '''4.513'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.21%  '
$ws.Range("D35").Value = '''2.669'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -5.51%  '
$ws.Range("E36").Value = '  +3.09%  '
$ws.Range("D37").Value = '''0.01976'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.14%  '
$ws.Range("E38").Value = '  +1.98%  '
$ws.Range("D39").Value = '''2.965'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.37%  '
$ws.Range("D40").Value = '''0.5287'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.25%  '
$ws.Range("D41").Value = '''7.033'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.03%  '
$ws.Range("D42").Value = '''0.1525'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.95%  '
$ws.Range("D43").Value = '''8.291'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.26%  '
$ws.Range("D44").Value = '''10.55'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.13%  '
$ws.Range("D45").Value = '''0.4743'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.36%  '
$ws.Range("E46").Value = '  +0.68%  '
$ws.Range("D47").Value = '''101.91'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.32%  '
$ws.Range("E48").Value = '  +0.74%  '
$ws.Range("D49").Value = '''66.05'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.49%  '
$ws.Range("D50").Value = '''0.06065'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.36%  '
$ws.Range("D51").Value = '''0.8942'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.89%  '
